# Scheduled-runner refresh of market-price derived columns (currentAveragePrice*,
# LevePrice*, LeveProfit*) across all 8 job sheets. Values below come from the
# latest Universalis price pull; only numeric cells move - no rows/cols/styles change.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 6426.7627
$ws.Range("I15").Value = 6426.7627
$ws.Range("K15").Value = 19280.2881
$ws.Range("M15").Value = -19111.2881

# row 28
$ws.Range("H28").Value = 676.3333
$ws.Range("I28").Value = 615.7143
$ws.Range("K28").Value = 615.7143
$ws.Range("M28").Value = -130.7143

# row 64
$ws.Range("H64").Value = 3396.5518
$ws.Range("I64").Value = 3275
$ws.Range("J64").Value = 3546.1538
$ws.Range("K64").Value = 3275
$ws.Range("L64").Value = 3546.1538
$ws.Range("M64").Value = -3027
$ws.Range("N64").Value = -4042.1538

# row 67
$ws.Range("H67").Value = 3396.5518
$ws.Range("I67").Value = 3275
$ws.Range("J67").Value = 3546.1538
$ws.Range("K67").Value = 3275
$ws.Range("L67").Value = 3546.1538
$ws.Range("M67").Value = -2417
$ws.Range("N67").Value = -5262.1538

# row 74
$ws.Range("H74").Value = 3691.9607
$ws.Range("I74").Value = 3632.162
$ws.Range("J74").Value = 3850
$ws.Range("K74").Value = 3632.162
$ws.Range("L74").Value = 3850
$ws.Range("M74").Value = -2696.162
$ws.Range("N74").Value = -5722

# row 76
$ws.Range("H76").Value = 3165.5293
$ws.Range("I76").Value = 3003
$ws.Range("J76").Value = 3187.2
$ws.Range("K76").Value = 3003
$ws.Range("L76").Value = 3187.2
$ws.Range("M76").Value = -2688
$ws.Range("N76").Value = -3817.2

# row 77
$ws.Range("H77").Value = 3691.9607
$ws.Range("I77").Value = 3632.162
$ws.Range("J77").Value = 3850
$ws.Range("K77").Value = 18160.81
$ws.Range("L77").Value = 19250
$ws.Range("M77").Value = -13480.81
$ws.Range("N77").Value = -28610

# row 79
$ws.Range("H79").Value = 3165.5293
$ws.Range("I79").Value = 3003
$ws.Range("J79").Value = 3187.2
$ws.Range("K79").Value = 3003
$ws.Range("L79").Value = 3187.2
$ws.Range("M79").Value = -1911
$ws.Range("N79").Value = -5371.2

# row 98
$ws.Range("H98").Value = 29999.703
$ws.Range("I98").Value = 711.3158
$ws.Range("K98").Value = 711.3158
$ws.Range("M98").Value = 786.6842

# row 122
$ws.Range("H122").Value = 29999.703
$ws.Range("I122").Value = 711.3158
$ws.Range("K122").Value = 2133.9474
$ws.Range("M122").Value = 316.0526

# row 132
$ws.Range("H132").Value = 21554.762
$ws.Range("I132").Value = 2655.9487
$ws.Range("K132").Value = 7967.8461
$ws.Range("M132").Value = -5437.8461

# row 137
$ws.Range("H137").Value = 1792422.2
$ws.Range("I137").Value = 4525714
$ws.Range("J137").Value = 5270.1924
$ws.Range("K137").Value = 13577142
$ws.Range("L137").Value = 15810.5772
$ws.Range("M137").Value = -13574592
$ws.Range("N137").Value = -20910.5772

# row 138
$ws.Range("H138").Value = 2712.4412
$ws.Range("I138").Value = 2544.2727
$ws.Range("J138").Value = 2792.8696
$ws.Range("K138").Value = 7632.8181
$ws.Range("L138").Value = 8378.6088
$ws.Range("M138").Value = -2492.8181
$ws.Range("N138").Value = -18658.6088

# row 141
$ws.Range("H141").Value = 3303.5881
$ws.Range("I141").Value = 1882.8889
$ws.Range("K141").Value = 5648.6667
$ws.Range("M141").Value = -468.6666999999998

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 1643.9524
$ws.Range("I45").Value = 1682.5
$ws.Range("J45").Value = 1608.909
$ws.Range("K45").Value = 1682.5
$ws.Range("L45").Value = 1608.909
$ws.Range("M45").Value = -1305.5
$ws.Range("N45").Value = -2362.909

# row 63
$ws.Range("H63").Value = 3300.625
$ws.Range("I63").Value = 2330
$ws.Range("K63").Value = 2330
$ws.Range("M63").Value = -1644

# row 66
$ws.Range("H66").Value = 3300.625
$ws.Range("I66").Value = 2330
$ws.Range("K66").Value = 11650
$ws.Range("M66").Value = -8218

# row 74
$ws.Range("H74").Value = 1721.0889
$ws.Range("I74").Value = 1343.0294
$ws.Range("K74").Value = 1343.0294
$ws.Range("M74").Value = -469.0293999999999

# row 77
$ws.Range("H77").Value = 1721.0889
$ws.Range("I77").Value = 1343.0294
$ws.Range("K77").Value = 6715.146999999999
$ws.Range("M77").Value = -2347.146999999999

# row 88
$ws.Range("H88").Value = 11616803
$ws.Range("I88").Value = 22225466
$ws.Range("J88").Value = 2936986.8
$ws.Range("K88").Value = 22225466
$ws.Range("L88").Value = 2936986.8
$ws.Range("M88").Value = -22225060
$ws.Range("N88").Value = -2937798.8

# row 91
$ws.Range("H91").Value = 11616803
$ws.Range("I91").Value = 22225466
$ws.Range("J91").Value = 2936986.8
$ws.Range("K91").Value = 22225466
$ws.Range("L91").Value = 2936986.8
$ws.Range("M91").Value = -22224062
$ws.Range("N91").Value = -2939794.8

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 2550
$ws.Range("I86").Value = 2500
$ws.Range("K86").Value = 2500
$ws.Range("M86").Value = -1377

# row 89
$ws.Range("H89").Value = 2550
$ws.Range("I89").Value = 2500
$ws.Range("K89").Value = 12500
$ws.Range("M89").Value = -6884

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# row 62
$ws.Range("H62").Value = 2800.5386
$ws.Range("J62").Value = 3121.2
$ws.Range("L62").Value = 3121.2
$ws.Range("N62").Value = -4369.2

# row 65
$ws.Range("H65").Value = 2800.5386
$ws.Range("J65").Value = 3121.2
$ws.Range("L65").Value = 15606
$ws.Range("N65").Value = -21846

# row 107
$ws.Range("H107").Value = 572.0789
$ws.Range("I107").Value = 450.86667
$ws.Range("J107").Value = 1026.625
$ws.Range("K107").Value = 450.86667
$ws.Range("L107").Value = 1026.625
$ws.Range("M107").Value = 1469.13333
$ws.Range("N107").Value = -4866.625

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# row 20
$ws.Range("H20").Value = 2850
$ws.Range("I20").Value = 466.66666
$ws.Range("K20").Value = 1399.99998
$ws.Range("M20").Value = -1172.99998

# row 107
$ws.Range("H107").Value = 347.73685
$ws.Range("J107").Value = 294.8
$ws.Range("L107").Value = 884.4000000000001
$ws.Range("N107").Value = -4724.4

# row 131
$ws.Range("H131").Value = 875.4693600000001

# row 132
$ws.Range("H132").Value = 2459.2778
$ws.Range("I132").Value = 965.375
$ws.Range("J132").Value = 3654.4
$ws.Range("K132").Value = 8688.375
$ws.Range("L132").Value = 32889.6
$ws.Range("M132").Value = -6158.375
$ws.Range("N132").Value = -37949.6

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 5508.846
$ws.Range("I70").Value = 5531.5
$ws.Range("J70").Value = 5433.3335
$ws.Range("K70").Value = 5531.5
$ws.Range("L70").Value = 5433.3335
$ws.Range("M70").Value = -5261.5
$ws.Range("N70").Value = -5973.3335

# row 73
$ws.Range("H73").Value = 5508.846
$ws.Range("I73").Value = 5531.5
$ws.Range("J73").Value = 5433.3335
$ws.Range("K73").Value = 5531.5
$ws.Range("L73").Value = 5433.3335
$ws.Range("M73").Value = -4595.5
$ws.Range("N73").Value = -7305.3335

# row 80
$ws.Range("H80").Value = 5485.7144
$ws.Range("I80").Value = 6142.857
$ws.Range("J80").Value = 4828.5713
$ws.Range("K80").Value = 6142.857
$ws.Range("L80").Value = 4828.5713
$ws.Range("M80").Value = -5144.857
$ws.Range("N80").Value = -6824.5713

# row 83
$ws.Range("H83").Value = 5485.7144
$ws.Range("I83").Value = 6142.857
$ws.Range("J83").Value = 4828.5713
$ws.Range("K83").Value = 30714.285
$ws.Range("L83").Value = 24142.8565
$ws.Range("M83").Value = -25722.285
$ws.Range("N83").Value = -34126.85649999999

# row 102
$ws.Range("H102").Value = 1000
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# row 61
$ws.Range("H61").Value = 2459.6667
$ws.Range("I61").Value = 2718.1
$ws.Range("J61").Value = 2224.7273
$ws.Range("K61").Value = 2718.1
$ws.Range("L61").Value = 2224.7273
$ws.Range("M61").Value = -2516.1
$ws.Range("N61").Value = -2628.7273

# row 113
$ws.Range("H113").Value = 2459.6667
$ws.Range("I113").Value = 2718.1
$ws.Range("J113").Value = 2224.7273
$ws.Range("K113").Value = 2718.1
$ws.Range("L113").Value = 2224.7273
$ws.Range("M113").Value = -548.0999999999999
$ws.Range("N113").Value = -6564.7273

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# row 113
$ws.Range("H113").Value = 390.9
$ws.Range("I113").Value = 367.66666
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 1102.99998
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 1067.00002
$ws.Range("N113").Value = -6140

